# Allow keys in permitted operations
# Replace the computed numeric results for rows whose "operation" value is
# not one of the canonically recognized operation names with the string
# "Not Implemented".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Not Implemented"
$ws.Range("D26").Value = "Not Implemented"
$ws.Range("D37").Value = "Not Implemented"
$ws.Range("D49").Value = "Not Implemented"
